$p = $ppt.ActivePresentation
$newDate = "2021/09/29"

function Set-DateFieldText($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$master = $p.SlideMaster
Set-DateFieldText $master

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Set-DateFieldText $layout
}

$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(7)
$shp.Left = 201.275390625
$shp.Top = 240.342845
